# Fill in the missing English (column E) translations for the Coroner-related
# localization rows that were added via the Crowdin translation pass.
# Column A = key, Column B = Japanese text (unchanged), Column E = English text (new).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$translations = @{
    26 = "Coroner";
    27 = "There is no ~r~dead bodies~s~ nearby you.";
    28 = "Requested ~b~{0}~s~ unit to Dispatch.";
    29 = "You can check ~b~Coroner's Report~s~ for more information.";
    30 = "Have a nice day! Officer!";
    31 = "Press {0} to teleport the backup unit nearby.";
    33 = "Coroner Menu";
    34 = "Coroner Report";
    35 = "Report Count: {0}";
    36 = "No Data";
    38 = "Name";
    39 = "Sex";
    40 = "Cause of Death";
    41 = "Died Day";
    43 = "Backup Vehicle";
    44 = "Backup Officer";
}

foreach ($row in ($translations.Keys | Sort-Object)) {
    $ws.Range("E$row").Value = $translations[$row]
}
